# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update status text everywhere it appears, so the shared string is updated in place
# (Overview sheet E2/F2/E3/F3 and the Status column (C2/C3) on the zh-cn/de-de sheets)
$wsOverview.Range("E2").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: not in sync with en-US"

$wsZhCn.Range("C2").Value = "Handed back: not in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: not in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: not in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: not in sync with en-US"

# Update Correspond Handback DateTime for ca05378c row (row 3) in zh-cn and de-de sheets
$wsZhCn.Range("K3").Value = "2016-11-01 14:37:11"
$wsDeDe.Range("K3").Value = "2016-11-01 14:37:28"
